$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.408218
$ws.Range("H2").Value = 61.224654
$ws.Range("I2").Value = 0.1108535210972707
$ws.Range("J2").Value = 0.1108535210972707
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 2973.492559919671
$ws.Range("R2").Value = 26761.43303927704
$ws.Range("S2").Value = 0.03176992908448226
$ws.Range("T2").Value = 0.03176992908448226
$ws.Range("G3").Value = 20.408218
$ws.Range("H3").Value = 61.224654
$ws.Range("I3").Value = 0.1108535210972707
$ws.Range("J3").Value = 0.1108535210972707
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 3444.901130356515
$ws.Range("R3").Value = 31004.11017320864
$ws.Range("S3").Value = 0.03680663812302792
$ws.Range("T3").Value = 0.03680663812302792
$ws.Range("G4").Value = 20.408218
$ws.Range("H4").Value = 61.224654
$ws.Range("I4").Value = 0.1108535210972707
$ws.Range("J4").Value = 0.1108535210972707
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 2614.826495939051
$ws.Range("R4").Value = 23533.43846345146
$ws.Range("S4").Value = 0.02793780400326715
$ws.Range("T4").Value = 0.02793780400326715
$ws.Range("G5").Value = 20.408218
$ws.Range("H5").Value = 61.224654
$ws.Range("I5").Value = 0.1108535210972707
$ws.Range("J5").Value = 0.1108535210972707
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 1342.066436147222
$ws.Range("R5").Value = 12078.597925325
$ws.Range("S5").Value = 0.01433914988649338
$ws.Range("T5").Value = 0.01433914988649338
$ws.Range("I6").Value = 0.2566851044076959
$ws.Range("J6").Value = 0.256685104407696
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 6885.223316711405
$ws.Range("R6").Value = 61967.00985040265
$ws.Range("S6").Value = 0.07356435306118754
$ws.Range("T6").Value = 0.07356435306118755
$ws.Range("I7").Value = 0.2566851044076959
$ws.Range("J7").Value = 0.256685104407696
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.08522702441914867
$ws.Range("T7").Value = 0.0852270244191487
$ws.Range("I8").Value = 0.2566851044076959
$ws.Range("J8").Value = 0.256685104407696
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 6054.719827340245
$ws.Range("R8").Value = 54492.4784460622
$ws.Range("S8").Value = 0.06469093689146636
$ws.Range("T8").Value = 0.06469093689146638
$ws.Range("I9").Value = 0.2566851044076959
$ws.Range("J9").Value = 0.256685104407696
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 3107.600551381996
$ws.Range("R9").Value = 27968.40496243796
$ws.Range("S9").Value = 0.03320279003589337
$ws.Range("T9").Value = 0.03320279003589338
$ws.Range("G10").Value = 85.307233
$ws.Range("H10").Value = 255.921699
$ws.Range("I10").Value = 0.4633725077375833
$ws.Range("J10").Value = 0.4633725077375833
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 12429.32737355284
$ws.Range("R10").Value = 111863.9463619756
$ws.Range("S10").Value = 0.1327996762286351
$ws.Range("T10").Value = 0.1327996762286352
$ws.Range("G11").Value = 85.307233
$ws.Range("H11").Value = 255.921699
$ws.Range("I11").Value = 0.4633725077375833
$ws.Range("J11").Value = 0.4633725077375833
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 14399.83556571606
$ws.Range("R11").Value = 129598.5200914445
$ws.Range("S11").Value = 0.1538533376264319
$ws.Range("T11").Value = 0.1538533376264319
$ws.Range("G12").Value = 85.307233
$ws.Range("H12").Value = 255.921699
$ws.Range("I12").Value = 0.4633725077375833
$ws.Range("J12").Value = 0.4633725077375833
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 10930.08772954337
$ws.Range("R12").Value = 98370.78956589033
$ws.Range("S12").Value = 0.1167812278178841
$ws.Range("T12").Value = 0.1167812278178841
$ws.Range("G13").Value = 85.307233
$ws.Range("H13").Value = 255.921699
$ws.Range("I13").Value = 0.4633725077375833
$ws.Range("J13").Value = 0.4633725077375833
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 5609.895688584406
$ws.Range("R13").Value = 50489.06119725966
$ws.Range("S13").Value = 0.05993826606463212
$ws.Range("T13").Value = 0.05993826606463213
$ws.Range("G14").Value = 31.12938966666666
$ws.Range("H14").Value = 93.38816899999999
$ws.Range("I14").Value = 0.16908886675745
$ws.Range("J14").Value = 0.16908886675745
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 4535.575255452172
$ws.Range("R14").Value = 40820.17729906955
$ws.Range("S14").Value = 0.04845981663627929
$ws.Range("T14").Value = 0.04845981663627929
$ws.Range("G15").Value = 31.12938966666666
$ws.Range("H15").Value = 93.38816899999999
$ws.Range("I15").Value = 0.16908886675745
$ws.Range("J15").Value = 0.16908886675745
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 5254.631719928139
$ws.Range("R15").Value = 47291.68547935325
$ws.Range("S15").Value = 0.0561424902679756
$ws.Range("T15").Value = 0.0561424902679756
$ws.Range("G16").Value = 31.12938966666666
$ws.Range("H16").Value = 93.38816899999999
$ws.Range("I16").Value = 0.16908886675745
$ws.Range("J16").Value = 0.16908886675745
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 3988.488995110269
$ws.Range("R16").Value = 35896.40095599242
$ws.Range("S16").Value = 0.04261453828299281
$ws.Range("T16").Value = 0.04261453828299281
$ws.Range("G17").Value = 31.12938966666666
$ws.Range("H17").Value = 93.38816899999999
$ws.Range("I17").Value = 0.16908886675745
$ws.Range("J17").Value = 0.16908886675745
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 2047.10225309145
$ws.Range("R17").Value = 18423.92027782305
$ws.Range("S17").Value = 0.02187202157020234
$ws.Range("T17").Value = 0.02187202157020234
